$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# NOTE: Excel's ColumnWidth (COM, character units) is offset from the stored
# OOXML <col width> by a constant ~0.83 (5-pixel default-font padding). Using
# "target - 0.83" reproduces the exact stored width.
$ws.Columns.Item(3).ColumnWidth = 63.17   # C: stored width 61 -> 64
$ws.Columns.Item(8).ColumnWidth = 51.17   # H: stored width 56 -> 52

# --- Ensure opportunity ID column (A) for all data rows stays text, not numeric ---
$idRange = $ws.Range("A2:A18")
$idRange.NumberFormat = "@"

# --- Update/insert row data (rows 2-18) ---
# Row 2: opportunity 1331343
$ws.Cells.Item(2, 1).Value = '1331343'
$ws.Cells.Item(2, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331343'
$ws.Cells.Item(2, 3).Value = 'Spanish Instructor for Children ( mexican or spanish only )'
$ws.Cells.Item(2, 4).Value = 'Londres, Royaume-Uni'
$ws.Cells.Item(2, 6).Value = '0 applicants'
$ws.Cells.Item(2, 7).Value = '6 - 18 Months'
$ws.Cells.Item(2, 8).Value = 'Bilingual Day Nursery and Preschool Ltd'

# Row 3: opportunity 1331341
$ws.Cells.Item(3, 1).Value = '1331341'
$ws.Cells.Item(3, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331341'
$ws.Cells.Item(3, 3).Value = 'Luxury Retail Intern'
$ws.Cells.Item(3, 4).Value = 'Londres, Royaume-Uni'
$ws.Cells.Item(3, 6).Value = '8 applicants'
$ws.Cells.Item(3, 7).Value = '3 - 6 Months'
$ws.Cells.Item(3, 8).Value = 'Liberty London'

# Row 4: opportunity 1331330
$ws.Cells.Item(4, 1).Value = '1331330'
$ws.Cells.Item(4, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331330'
$ws.Cells.Item(4, 3).Value = 'Accelerate Serbia | Software developer intern'
$ws.Cells.Item(4, 4).Value = 'Београд, Србија'
$ws.Cells.Item(4, 6).Value = '1 applicant'
$ws.Cells.Item(4, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(4, 8).Value = 'Privredni savetnik - Informacione tehnologije doo'

# Row 5: opportunity 1331322
$ws.Cells.Item(5, 1).Value = '1331322'
$ws.Cells.Item(5, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331322'
$ws.Cells.Item(5, 3).Value = 'Social Media Marketing'
$ws.Cells.Item(5, 4).Value = 'Petaling Jaya, Selangor, Malaysia'
$ws.Cells.Item(5, 6).Value = '1 applicant'
$ws.Cells.Item(5, 7).Value = '6 - 18 Months'
$ws.Cells.Item(5, 8).Value = 'iWisers SDN BHD'

# Row 6: opportunity 1331321
$ws.Cells.Item(6, 1).Value = '1331321'
$ws.Cells.Item(6, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331321'
$ws.Cells.Item(6, 3).Value = '[DGF] Global Marketing Intern - Content & Campaign Management'
$ws.Cells.Item(6, 4).Value = 'Charles-de-Gaulle-Straße 20, 53113 Bonn, Germany'
$ws.Cells.Item(6, 6).Value = '24 applicants'
$ws.Cells.Item(6, 7).Value = '6 - 18 Months'
$ws.Cells.Item(6, 8).Value = 'DHL Group'

# Row 7: opportunity 1331270
$ws.Cells.Item(7, 1).Value = '1331270'
$ws.Cells.Item(7, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331270'
$ws.Cells.Item(7, 3).Value = 'Digital Marketing Intern'
$ws.Cells.Item(7, 4).Value = 'Sahibzada Ajit Singh Nagar, Punjab, India'
$ws.Cells.Item(7, 6).Value = '0 applicants'
$ws.Cells.Item(7, 7).Value = '3 - 6 Months'
$ws.Cells.Item(7, 8).Value = 'Requisite Technologies Pvt Ltd'

# Row 8: opportunity 1331071
$ws.Cells.Item(8, 1).Value = '1331071'
$ws.Cells.Item(8, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331071'
$ws.Cells.Item(8, 3).Value = 'Growth Analytics & Funnel Optimization'
$ws.Cells.Item(8, 4).Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Cells.Item(8, 6).Value = '0 applicants'
$ws.Cells.Item(8, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(8, 8).Value = 'Madaar'

# Row 9: opportunity 1331070
$ws.Cells.Item(9, 1).Value = '1331070'
$ws.Cells.Item(9, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331070'
$ws.Cells.Item(9, 3).Value = 'Performance Marketing'
$ws.Cells.Item(9, 4).Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Cells.Item(9, 6).Value = '0 applicants'
$ws.Cells.Item(9, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(9, 8).Value = 'Madaar'

# Row 10: opportunity 1331069
$ws.Cells.Item(10, 1).Value = '1331069'
$ws.Cells.Item(10, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331069'
$ws.Cells.Item(10, 3).Value = 'Growth Marketing'
$ws.Cells.Item(10, 4).Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Cells.Item(10, 6).Value = '0 applicants'
$ws.Cells.Item(10, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(10, 8).Value = 'Madaar'

# Row 11: opportunity 1331068
$ws.Cells.Item(11, 1).Value = '1331068'
$ws.Cells.Item(11, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331068'
$ws.Cells.Item(11, 3).Value = 'Product Marketing'
$ws.Cells.Item(11, 4).Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Cells.Item(11, 6).Value = '0 applicants'
$ws.Cells.Item(11, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(11, 8).Value = 'Madaar'

# Row 12: opportunity 1331067
$ws.Cells.Item(12, 1).Value = '1331067'
$ws.Cells.Item(12, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331067'
$ws.Cells.Item(12, 3).Value = 'Business Development'
$ws.Cells.Item(12, 4).Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Cells.Item(12, 6).Value = '0 applicants'
$ws.Cells.Item(12, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(12, 8).Value = 'Madaar'

# Row 13: opportunity 1331066
$ws.Cells.Item(13, 1).Value = '1331066'
$ws.Cells.Item(13, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331066'
$ws.Cells.Item(13, 3).Value = 'Sales'
$ws.Cells.Item(13, 4).Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Cells.Item(13, 6).Value = '0 applicants'
$ws.Cells.Item(13, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(13, 8).Value = 'Madaar'

# Row 14: opportunity 1331065
$ws.Cells.Item(14, 1).Value = '1331065'
$ws.Cells.Item(14, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331065'
$ws.Cells.Item(14, 3).Value = 'Content & Brand Marketing'
$ws.Cells.Item(14, 4).Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Cells.Item(14, 6).Value = '0 applicants'
$ws.Cells.Item(14, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(14, 8).Value = 'Madaar'

# Row 15: opportunity 1331064
$ws.Cells.Item(15, 1).Value = '1331064'
$ws.Cells.Item(15, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331064'
$ws.Cells.Item(15, 3).Value = 'Video Editor / Reels Maker'
$ws.Cells.Item(15, 4).Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Cells.Item(15, 6).Value = '1 applicant'
$ws.Cells.Item(15, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(15, 8).Value = 'Madaar'

# Row 16: opportunity 1331063
$ws.Cells.Item(16, 1).Value = '1331063'
$ws.Cells.Item(16, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331063'
$ws.Cells.Item(16, 3).Value = 'Social Media & Community Marketing'
$ws.Cells.Item(16, 4).Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Cells.Item(16, 6).Value = '0 applicants'
$ws.Cells.Item(16, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(16, 8).Value = 'Madaar'

# Row 17: opportunity 1328549
$ws.Cells.Item(17, 1).Value = '1328549'
$ws.Cells.Item(17, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328549'
$ws.Cells.Item(17, 3).Value = 'Motion graphic design'
$ws.Cells.Item(17, 4).Value = 'Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt'
$ws.Cells.Item(17, 6).Value = '0 applicants'
$ws.Cells.Item(17, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(17, 8).Value = 'TAR - Company'

# Row 18: opportunity 1325831
$ws.Cells.Item(18, 1).Value = '1325831'
$ws.Cells.Item(18, 2).Value = 'https://aiesec.org/opportunity/global-talent/1325831'
$ws.Cells.Item(18, 3).Value = 'Information Technology Intern'
$ws.Cells.Item(18, 4).Value = 'Makati City, Metro Manila, Philippines'
$ws.Cells.Item(18, 6).Value = '30 applicants'
$ws.Cells.Item(18, 7).Value = '3 - 6 Months'
$ws.Cells.Item(18, 8).Value = 'Consistent Frozen Solutions Corporation'

# --- Restore column A default style (drop the temporary text-number format) ---
$idRange.Style = "Normal"

# --- PREMIUM column (E): set value + fill matching "Yes" (yellow highlight) / "No" (plain) ---
# Use existing styled cells as format templates so no duplicate styles get created.

$ws.Range("E4").Copy($ws.Cells.Item(2, 5))   # row 2: PREMIUM = No (plain style)
$ws.Cells.Item(2, 5).Value = 'No'
$ws.Range("E4").Copy($ws.Cells.Item(3, 5))   # row 3: PREMIUM = No (plain style)
$ws.Cells.Item(3, 5).Value = 'No'
$ws.Range("E4").Copy($ws.Cells.Item(4, 5))   # row 4: PREMIUM = No (plain style)
$ws.Cells.Item(4, 5).Value = 'No'
$ws.Range("E4").Copy($ws.Cells.Item(5, 5))   # row 5: PREMIUM = No (plain style)
$ws.Cells.Item(5, 5).Value = 'No'
$ws.Range("E2").Copy($ws.Cells.Item(6, 5))   # row 6: PREMIUM = Yes (apply yellow highlight style)
$ws.Cells.Item(6, 5).Value = 'Yes'
$ws.Range("E4").Copy($ws.Cells.Item(7, 5))   # row 7: PREMIUM = No (plain style)
$ws.Cells.Item(7, 5).Value = 'No'
$ws.Range("E4").Copy($ws.Cells.Item(8, 5))   # row 8: PREMIUM = No (plain style)
$ws.Cells.Item(8, 5).Value = 'No'
$ws.Range("E4").Copy($ws.Cells.Item(9, 5))   # row 9: PREMIUM = No (plain style)
$ws.Cells.Item(9, 5).Value = 'No'
$ws.Range("E4").Copy($ws.Cells.Item(10, 5))   # row 10: PREMIUM = No (plain style)
$ws.Cells.Item(10, 5).Value = 'No'
$ws.Range("E4").Copy($ws.Cells.Item(11, 5))   # row 11: PREMIUM = No (plain style)
$ws.Cells.Item(11, 5).Value = 'No'
$ws.Range("E4").Copy($ws.Cells.Item(12, 5))   # row 12: PREMIUM = No (plain style)
$ws.Cells.Item(12, 5).Value = 'No'
$ws.Range("E4").Copy($ws.Cells.Item(13, 5))   # row 13: PREMIUM = No (plain style)
$ws.Cells.Item(13, 5).Value = 'No'
$ws.Range("E4").Copy($ws.Cells.Item(14, 5))   # row 14: PREMIUM = No (plain style)
$ws.Cells.Item(14, 5).Value = 'No'
$ws.Range("E4").Copy($ws.Cells.Item(15, 5))   # row 15: PREMIUM = No (plain style)
$ws.Cells.Item(15, 5).Value = 'No'
$ws.Range("E4").Copy($ws.Cells.Item(16, 5))   # row 16: PREMIUM = No (plain style)
$ws.Cells.Item(16, 5).Value = 'No'
$ws.Range("E4").Copy($ws.Cells.Item(17, 5))   # row 17: PREMIUM = No (plain style)
$ws.Cells.Item(17, 5).Value = 'No'
$ws.Range("E4").Copy($ws.Cells.Item(18, 5))   # row 18: PREMIUM = No (plain style)
$ws.Cells.Item(18, 5).Value = 'No'

$excel.CutCopyMode = $false

Write-Host "Edit complete"
